# "Energy Single Operation" column: switch the unit from Joule to picojoule.
#  - header label [J] -> [pJ]
#  - formula simplified: the old formula pre-converted to Joules
#    ((K/1000)*X*(10^-9)); the new one leaves the result in picojoules
#    (K*X), since pJ = (J) * 1e12 and the old formula's result *1e12
#    collapses the (10^-9)/1000 factor away, i.e. M = K*X.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Header text.
$ws.Range("M3").Value = "Energy Single Operation [pJ]"

# Re-point every row's formula from the old J-with-prefactor expression
# to the simplified pJ one ($K*$X, same result scaled by 1e12).
for ($r = 5; $r -le 22; $r++) {
    $ws.Range("M$r").Formula = "=`$K$r*`$X$r"
}

# Restore the selection left by the editor (was I25, now N30).
$ws.Activate() | Out-Null
$ws.Range("N30").Select() | Out-Null
